# chore: update Sheets via scheduled runner
# Refresh cached market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 100.63636
$ws.Range("I33").Value = 89.666664
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 89.666664
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 139.333336
$ws.Range("N33").Value = -608

$ws.Range("H40").Value = 2793
$ws.Range("J40").Value = 2991
$ws.Range("L40").Value = 2991
$ws.Range("N40").Value = -3341

$ws.Range("H64").Value = 62503012
$ws.Range("I64").Value = 166668860
$ws.Range("K64").Value = 166668860
$ws.Range("M64").Value = -166668612

$ws.Range("H67").Value = 62503012
$ws.Range("I67").Value = 166668860
$ws.Range("K67").Value = 166668860
$ws.Range("M67").Value = -166668002

$ws.Range("H69").Value = 3903
$ws.Range("I69").Value = 3910
$ws.Range("J69").Value = 3899.5
$ws.Range("K69").Value = 11730
$ws.Range("L69").Value = 11698.5
$ws.Range("M69").Value = -10856
$ws.Range("N69").Value = -13446.5

$ws.Range("H70").Value = 7875.5
$ws.Range("J70").Value = 4250
$ws.Range("L70").Value = 12750
$ws.Range("N70").Value = -13290

$ws.Range("H72").Value = 3903
$ws.Range("I72").Value = 3910
$ws.Range("J72").Value = 3899.5
$ws.Range("K72").Value = 35190
$ws.Range("L72").Value = 35095.5
$ws.Range("M72").Value = -30822
$ws.Range("N72").Value = -43831.5

$ws.Range("H73").Value = 7875.5
$ws.Range("J73").Value = 4250
$ws.Range("L73").Value = 12750
$ws.Range("N73").Value = -14622

$ws.Range("H74").Value = 3131
$ws.Range("I74").Value = 2430.8
$ws.Range("J74").Value = 3449.2727
$ws.Range("K74").Value = 2430.8
$ws.Range("L74").Value = 3449.2727
$ws.Range("M74").Value = -1494.8
$ws.Range("N74").Value = -5321.2727

$ws.Range("H76").Value = 5939
$ws.Range("I76").Value = 4167
$ws.Range("J76").Value = 6293.4
$ws.Range("K76").Value = 4167
$ws.Range("L76").Value = 6293.4
$ws.Range("M76").Value = -3852
$ws.Range("N76").Value = -6923.4

$ws.Range("H77").Value = 3131
$ws.Range("I77").Value = 2430.8
$ws.Range("J77").Value = 3449.2727
$ws.Range("K77").Value = 12154
$ws.Range("L77").Value = 17246.3635
$ws.Range("M77").Value = -7474
$ws.Range("N77").Value = -26606.3635

$ws.Range("H79").Value = 5939
$ws.Range("I79").Value = 4167
$ws.Range("J79").Value = 6293.4
$ws.Range("K79").Value = 4167
$ws.Range("L79").Value = 6293.4
$ws.Range("M79").Value = -3075
$ws.Range("N79").Value = -8477.4

$ws.Range("H81").Value = 30001
$ws.Range("J81").Value = 30001
$ws.Range("L81").Value = 30001
$ws.Range("N81").Value = -31997

$ws.Range("H84").Value = 30001
$ws.Range("J84").Value = 30001
$ws.Range("L84").Value = 90003
$ws.Range("N84").Value = -99987

$ws.Range("H132").Value = 2937.673
$ws.Range("I132").Value = 2388.4866
$ws.Range("J132").Value = 4292.3335
$ws.Range("K132").Value = 7165.459800000001
$ws.Range("L132").Value = 12877.0005
$ws.Range("M132").Value = -4635.459800000001
$ws.Range("N132").Value = -17937.0005

$ws.Range("H137").Value = 3339.6826
$ws.Range("I137").Value = 1548
$ws.Range("J137").Value = 4056.3555
$ws.Range("K137").Value = 4644
$ws.Range("L137").Value = 12169.0665
$ws.Range("M137").Value = -2094
$ws.Range("N137").Value = -17269.0665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10412.884
$ws.Range("I32").Value = 8777.641
$ws.Range("J32").Value = 34669
$ws.Range("K32").Value = 8777.641
$ws.Range("L32").Value = 34669
$ws.Range("M32").Value = -8490.641
$ws.Range("N32").Value = -35243

$ws.Range("H74").Value = 3122.0425
$ws.Range("I74").Value = 892.35
$ws.Range("J74").Value = 15863.143
$ws.Range("K74").Value = 892.35
$ws.Range("L74").Value = 15863.143
$ws.Range("M74").Value = -18.35000000000002
$ws.Range("N74").Value = -17611.143

$ws.Range("H77").Value = 3122.0425
$ws.Range("I77").Value = 892.35
$ws.Range("J77").Value = 15863.143
$ws.Range("K77").Value = 4461.75
$ws.Range("L77").Value = 79315.715
$ws.Range("M77").Value = -93.75
$ws.Range("N77").Value = -88051.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 10000
$ws.Range("K34").Value = 10000
$ws.Range("M34").Value = -9886

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 181.84
$ws.Range("I7").Value = 255.64285
$ws.Range("J7").Value = 87.90909000000001
$ws.Range("K7").Value = 255.64285
$ws.Range("L7").Value = 87.90909000000001
$ws.Range("M7").Value = -142.64285
$ws.Range("N7").Value = -313.90909

$ws.Range("H31").Value = 14305066
$ws.Range("I31").Value = 66668020
$ws.Range("J31").Value = 24260.908
$ws.Range("K31").Value = 66668020
$ws.Range("L31").Value = 24260.908
$ws.Range("M31").Value = -66667725
$ws.Range("N31").Value = -24850.908

$ws.Range("H34").Value = 14305066
$ws.Range("I34").Value = 66668020
$ws.Range("J34").Value = 24260.908
$ws.Range("K34").Value = 66668020
$ws.Range("L34").Value = 24260.908
$ws.Range("M34").Value = -66667818
$ws.Range("N34").Value = -24664.908

$ws.Range("H62").Value = 12390
$ws.Range("I62").Value = 16200
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 16200
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -15576
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 12390
$ws.Range("I65").Value = 16200
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 81000
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -77880
$ws.Range("N65").Value = -23740

$ws.Range("H132").Value = 5319.275
$ws.Range("I132").Value = 5856.913
$ws.Range("J132").Value = 4591.8823
$ws.Range("K132").Value = 17570.739
$ws.Range("L132").Value = 13775.6469
$ws.Range("M132").Value = -15040.739
$ws.Range("N132").Value = -18835.6469

$ws.Range("H141").Value = 41699840
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 47655532
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 47655532
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -47665892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3562.9666
$ws.Range("I70").Value = 3294.9546
$ws.Range("J70").Value = 4300
$ws.Range("K70").Value = 3294.9546
$ws.Range("L70").Value = 4300
$ws.Range("M70").Value = -3024.9546
$ws.Range("N70").Value = -4840

$ws.Range("H73").Value = 3562.9666
$ws.Range("I73").Value = 3294.9546
$ws.Range("J73").Value = 4300
$ws.Range("K73").Value = 3294.9546
$ws.Range("L73").Value = 4300
$ws.Range("M73").Value = -2358.9546
$ws.Range("N73").Value = -6172

$ws.Range("H134").Value = 21046
$ws.Range("J134").Value = 21046
$ws.Range("L134").Value = 63138
$ws.Range("N134").Value = -68208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1496.25
$ws.Range("I7").Value = 1093
$ws.Range("J7").Value = 2168.3333
$ws.Range("K7").Value = 1093
$ws.Range("L7").Value = 2168.3333
$ws.Range("M7").Value = -981
$ws.Range("N7").Value = -2392.3333

$ws.Range("H68").Value = 1811.862
$ws.Range("I68").Value = 1550.7333
$ws.Range("J68").Value = 2091.6428
$ws.Range("K68").Value = 1550.7333
$ws.Range("L68").Value = 2091.6428
$ws.Range("M68").Value = -801.7333000000001
$ws.Range("N68").Value = -3589.6428

$ws.Range("H71").Value = 1811.862
$ws.Range("I71").Value = 1550.7333
$ws.Range("J71").Value = 2091.6428
$ws.Range("K71").Value = 7753.6665
$ws.Range("L71").Value = 10458.214
$ws.Range("M71").Value = -4009.6665
$ws.Range("N71").Value = -17946.214

$ws.Range("H126").Value = 1496.25
$ws.Range("I126").Value = 1093
$ws.Range("J126").Value = 2168.3333
$ws.Range("K126").Value = 3279
$ws.Range("L126").Value = 6504.999899999999
$ws.Range("M126").Value = -809
$ws.Range("N126").Value = -11444.9999

$ws.Range("H136").Value = 6299.857
$ws.Range("I136").Value = 2531.0557
$ws.Range("J136").Value = 13083.7
$ws.Range("K136").Value = 7593.1671
$ws.Range("L136").Value = 39251.10000000001
$ws.Range("M136").Value = -5043.1671
$ws.Range("N136").Value = -44351.10000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 465.53333
$ws.Range("I126").Value = 292.86957
$ws.Range("K126").Value = 878.60871
$ws.Range("M126").Value = 1591.39129
